$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is updated with the latest daily spot-price data (automatic price refresh)

$ws.Range("A2").Value = 45952

$ws.Range("B2").Value  = 37.6
$ws.Range("C2").Value  = 30.54
$ws.Range("D2").Value  = 27.62
$ws.Range("E2").Value  = 24.3
$ws.Range("F2").Value  = 24.28
$ws.Range("G2").Value  = 27.58
$ws.Range("H2").Value  = 43.05
$ws.Range("I2").Value  = 56.32
$ws.Range("J2").Value  = 74.39
$ws.Range("K2").Value  = 46.63
$ws.Range("L2").Value  = 23.29
$ws.Range("M2").Value  = 3.39
$ws.Range("N2").Value  = 1.02
$ws.Range("O2").Value  = 0.39
$ws.Range("P2").Value  = 0.03
$ws.Range("Q2").Value  = 0.01
$ws.Range("R2").Value  = 0.38
$ws.Range("S2").Value  = 3.43
$ws.Range("T2").Value  = 31.2
$ws.Range("U2").Value  = 60.14
$ws.Range("V2").Value  = 66.61
$ws.Range("W2").Value  = 43.19
$ws.Range("X2").Value  = 26.26
$ws.Range("Y2").Value  = 20
$ws.Range("Z2").Value  = 27.99

# AA2 (Slot_4h_max label) is unchanged ("20h-24h")

$ws.Range("AB2").Value = 39.02

$ws.Range("AC2").Value = "8h-10h"
$ws.Range("AD2").Value = 60.51
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 54.9
$ws.Range("AG2").Value = "2h-23h"
